# "question 2 in 12th" - automate/12th/res/example.xlsx
#
# 1. Insert two new blank rows above the "Pears" row (old row 3), pushing
#    the remaining data down two rows.
# 2. Strip all formatting from the last two data rows (they keep their
#    values, but revert to the workbook's default/no style).
# 3. Apply a custom date/time number format to the Date column for the
#    rows that used to sit below the inserted rows (old rows 3-7, i.e.
#    now rows 5-9) - this recreates styles for cells that had formatting,
#    and a fresh minimal style for the previously-cleared A8:A9 cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 2 blank rows before row 3 (shifts old rows 3-7 down to rows 5-9)
$ws.Rows("3:4").Insert()

# Remove formatting from the last two data rows entirely (Date/Item/Count)
$ws.Range("A8:C9").ClearFormats()

# Apply a custom date-time number format to the Date column (A) of the
# rows that used to be rows 3-7 (now rows 5-9)
$ws.Range("A5:A9").NumberFormat = "yyyy-mm-dd h:mm:ss"
